$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: add one "공정방식(B80)" common-code row, mirroring the format already
# used by the existing rows (REG_ID/REG_DTM columns re-use the style applied
# to the prior last row, 124).
function Add-CodeRow($Row, $A, $B, $C, $D, $E) {
    $ws.Rows.Item($Row).RowHeight = 15.75

    $ws.Cells.Item($Row, 1).Value = $A
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C

    $ws.Cells.Item(124, 4).Copy()
    $ws.Cells.Item($Row, 4).PasteSpecial(-4122)
    $ws.Cells.Item($Row, 4).Value = $D

    $ws.Cells.Item(124, 5).Copy()
    $ws.Cells.Item($Row, 5).PasteSpecial(-4122)
    $ws.Cells.Item($Row, 5).Value = $E
}

Add-CodeRow 125 "B80" "`$`$" "전체공정" "admin" 20201201153328
Add-CodeRow 126 "B80" "10" "전처리공정" "admin" 20201201153329
Add-CodeRow 127 "B80" "15" "용해공정" "admin" 20201201153330
Add-CodeRow 128 "B80" "20" "혼합공정" "admin" 20201201153331
Add-CodeRow 129 "B80" "25" "충진공정" "admin" 20201201153332
Add-CodeRow 130 "B80" "30" "가열공정" "admin" 20201201153333
Add-CodeRow 131 "B80" "35" "덧바르기" "admin" 20201201153334
Add-CodeRow 132 "B80" "36" "덧가루공정" "admin" 20201201153335
Add-CodeRow 133 "B80" "37" "반죽물공정" "admin" 20201201153336
Add-CodeRow 134 "B80" "38" "빵가루공정" "admin" 20201201153337
Add-CodeRow 135 "B80" "40" "살균공정" "admin" 20201201153338
Add-CodeRow 136 "B80" "45" "세척공정" "admin" 20201201153339
Add-CodeRow 137 "B80" "50" "냉장보관" "admin" 20201201153340
Add-CodeRow 138 "B80" "55" "토핑공정" "admin" 20201201153341
Add-CodeRow 139 "B80" "60" "절단공정" "admin" 20201201153342
Add-CodeRow 140 "B80" "65" "탈수공정" "admin" 20201201153343
Add-CodeRow 141 "B80" "70" "취반공정" "admin" 20201201153344
Add-CodeRow 142 "B80" "75" "성형공정" "admin" 20201201153345
Add-CodeRow 143 "B80" "80" "식힘공정" "admin" 20201201153346
Add-CodeRow 144 "B80" "85" "포장공정" "admin" 20201201153347
Add-CodeRow 145 "B80" "90" "금속검출공정" "admin" 20201201153348
Add-CodeRow 146 "B80" "95" "급냉공정" "admin" 20201201153349
Add-CodeRow 147 "B80" "99" "출하공정" "admin" 20201201153350

# Match the author's final selection (view scrolled down to the newly added data).
$ws.Range("B144").Select() | Out-Null
